# Report - data support.xlsx edit
# 1. Fill in bank-name / city detail columns (H, I) on the "Europe" sheet
#    for BIC rows that were missing them.
# 2. Update the remembered selection on "sSanctioned" (G4 -> J8).
# 3. Update the remembered selection on "Europe" (H121 -> H128) and make
#    "Europe" the active sheet/tab (was "Summary").

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Europe")

$ws.Range("H84").Value = "COMMERZBANK AG "
$ws.Range("I84").Value = "PARIS"
$ws.Range("H85").Value = "COMMERZBANK AG "
$ws.Range("I85").Value = "LUXEMBOURG"
$ws.Range("H86").Value = "SPARKASSE KOELNBONN"
$ws.Range("I86").Value = "KOELN"
$ws.Range("H87").Value = "BANK OF COMMUNICATIONS CO. LTD"
$ws.Range("I87").Value = "FRANKFURT "
$ws.Range("H88").Value = "CREDIT SUISSE (SCHWEIZ) AG"
$ws.Range("I88").Value = "ZURICH"
$ws.Range("H89").Value = "DEUTSCHE BANK (SUISSE) SA – PRIVATE BANKING"
$ws.Range("I89").Value = "GENEVA"
$ws.Range("H90").Value = "DEUTSCHE BANK AG"
$ws.Range("I90").Value = "FRANKFURT "
$ws.Range("H91").Value = "COMMERZBANK AG (FORMERLY DRESDNER BANK AG)"
$ws.Range("I91").Value = "FRANKFURT "
$ws.Range("H92").Value = "STADTSPARKASSE DUESSELDORF"
$ws.Range("I92").Value = "DUESSELDORF"
$ws.Range("H93").Value = "EBI SA "
$ws.Range("I93").Value = "PARIS"
$ws.Range("H94").Value = "EUROPEAN FUND ADMINISTRATION S.A."
$ws.Range("I94").Value = "LUXEMBOURG"
$ws.Range("H95").Value = "EUROPAEISCH-IRANISCHE HANDELSBANK AG"
$ws.Range("I95").Value = "HAMBURG"
$ws.Range("H97").Value = "RBC INVESTOR SERVICES BANK S.A."
$ws.Range("I97").Value = "ESCH-SUR-ALZETTE"
$ws.Range("H98").Value = "FIBI BANK (SWITZERLAND) LTD."
$ws.Range("I98").Value = "ZURICH"
$ws.Range("H99").Value = "CACEIS BANK"
$ws.Range("I99").Value = "LUXEMBOURG"
$ws.Range("H100").Value = "FRANKLIN TEMPLETON INTERNATIONAL SERVICES S.A."
$ws.Range("I100").Value = "LUXEMBOURG"
$ws.Range("H101").Value = "DZ BANK AG (FORMERLY WGZ BANK AG)"
$ws.Range("I101").Value = "DUESSELDORF"
$ws.Range("H102").Value = "DZ BANK AG, DEUTSCHE ZENTRAL- GENOSSENSCHAFTSBANK"
$ws.Range("I102").Value = "FRANKFURT "
$ws.Range("H103").Value = "BELFIUS BANK SA/NV"
$ws.Range("I103").Value = "BRUSSELS"
$ws.Range("H104").Value = "GUTZWILLER E. ET CIE BANQUIERS"
$ws.Range("I104").Value = "BASLE"
$ws.Range("H105").Value = "HABIB BANK AG"
$ws.Range("I105").Value = "ZURICH"
$ws.Range("H106").Value = "BANQUE HERITAGE"
$ws.Range("I106").Value = "GENEVA"
$ws.Range("H107").Value = "HSH NORDBANK AG"
$ws.Range("I107").Value = "HAMBURG"
$ws.Range("H108").Value = "UNICREDIT BANK AG (HYPOVEREINSBANK)"
$ws.Range("I108").Value = "MUENCHEN"
$ws.Range("H109").Value = "ICICI BANK UK PLC"
$ws.Range("I109").Value = "FRANKFURT "
$ws.Range("H110").Value = "INCORE BANK AG"
$ws.Range("I110").Value = "ZURICH"
$ws.Range("H111").Value = "SIX SIS AG"
$ws.Range("I111").Value = "ZURICH"
$ws.Range("H112").Value = "THE BANK OF NEW YORK MELLON (LUXEMBOURG) SA"
$ws.Range("I112").Value = "LUXEMBOURG"
$ws.Range("H113").Value = "BERNER KANTONALBANK AG"
$ws.Range("I113").Value = "BERNE"
$ws.Range("H114").Value = "KBL EUROPEAN PRIVATE BANKERS S.A."
$ws.Range("I114").Value = "LUXEMBOURG"
$ws.Range("H115").Value = "KBC SECURITIES N.V."
$ws.Range("I115").Value = "BRUSSELS"
$ws.Range("H116").Value = "ZUGER KANTONALBANK"
$ws.Range("I116").Value = "ZUG"
$ws.Range("H117").Value = "KBC BANK NV"
$ws.Range("I117").Value = "ANTWERPEN"
$ws.Range("H118").Value = "KBC BANK NV"
$ws.Range("I118").Value = "BRUSSELS"
$ws.Range("H119").Value = "KBL (SWITZERLAND) LTD"
$ws.Range("I119").Value = "GENEVA"
$ws.Range("H120").Value = "BANKHAUS LAMPE KG"
$ws.Range("I120").Value = "DUESSELDORF "
$ws.Range("H122").Value = "BANQUE LOMBARD ODIER ET CIE SA"
$ws.Range("I122").Value = "GENEVA"
$ws.Range("H123").Value = "LAROCHE ET CO. BANQUIERS"
$ws.Range("I123").Value = "BASLE"
$ws.Range("H124").Value = "LUZERNER KANTONALBANK"
$ws.Range("I124").Value = "LUCERNE"
$ws.Range("H125").Value = "BANK LEUMI (SWITZERLAND) LTD"
$ws.Range("I125").Value = "ZURICH"
$ws.Range("H126").Value = "MAN SE"
$ws.Range("I126").Value = "MUENCHEN "
$ws.Range("H127").Value = "BANKMED (SUISSE) S.A."
$ws.Range("I127").Value = "GENEVA"

# --- Selections / active sheet ---

# sSanctioned: selection moves from G4 to J8
$wsSanctioned = $wb.Worksheets.Item("sSanctioned")
$wsSanctioned.Activate()
$wsSanctioned.Range("J8").Select()

# Europe: selection moves from H121 to H128, and Europe becomes the
# active/selected tab (previously "Summary" was active).
$ws.Activate()
$ws.Range("H128").Select()
